$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.967.86'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.051.59'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '518.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.56%  '
$ws.Range("E8").Value = '  +3.94%  '
$ws.Range("E9").Value = '  +2.97%  '
$ws.Range("E10").Value = '  +5.49%  '
$ws.Range("E11").Value = '  +5.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.577.55'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.53%  '
$ws.Range("E13").Value = '  +2.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.76'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.71%  '
$ws.Range("E15").Value = '  +12.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '57.976.49'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.78%  '
$ws.Range("E17").Value = '  +9.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.050.78'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.05'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.89%  '
$ws.Range("E20").Value = '  +4.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '338.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.50%  '
$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.78'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.52%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("E24").Value = '  +7.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.09'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.21%  '
$ws.Range("E26").Value = '  +3.89%  '
$ws.Range("E27").Value = '  +7.74%  '
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("E29").Value = '  +6.46%  '
$ws.Range("E30").Value = '  +12.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.87%  '
$ws.Range("E32").Value = '  +4.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.04'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.14%  '
$ws.Range("E34").Value = '  +7.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.72'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.92'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.23'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.18%  '
$ws.Range("E39").Value = '  +2.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.087.18'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.80'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.15%  '
$ws.Range("E42").Value = '  +10.10%  '
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("E44").Value = '  +3.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.328.98'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.77%  '
$ws.Range("E46").Value = '  +4.40%  '
$ws.Range("E47").Value = '  +2.48%  '
$ws.Range("E48").Value = '  +5.46%  '
$ws.Range("E49").Value = '  +3.18%  '
$ws.Range("E50").Value = '  +4.89%  '
$ws.Range("E51").Value = '  -2.91%  '
